$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that is the same for every
# data row (46060 = 2026-02-07). This update bumps it by one day
# (-> 46061 = 2026-02-08) for every row, reflecting a refreshed
# "last changed" timestamp.
$lastRow = $ws.Cells.Item(1, 3).End(-4121).Row  # xlDown = -4121

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value = 46061
    }
}
